{"js": "// Update each two-digit multiplication problem's operands in the table.\n// Every \"A\u00d7B=\" string occurs exactly once in the document, so a scoped\n// search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"15\u00d756=\", \"75\u00d728=\"],\n  [\"74\u00d737=\", \"53\u00d716=\"],\n  [\"69\u00d721=\", \"41\u00d771=\"],\n  [\"19\u00d717=\", \"93\u00d787=\"],\n  [\"30\u00d778=\", \"55\u00d757=\"],\n  [\"31\u00d791=\", \"91\u00d768=\"],\n  [\"40\u00d756=\", \"20\u00d783=\"],\n  [\"64\u00d753=\", \"57\u00d796=\"],\n  [\"76\u00d712=\", \"97\u00d716=\"],\n  [\"50\u00d741=\", \"35\u00d786=\"],\n  [\"22\u00d732=\", \"89\u00d789=\"],\n  [\"91\u00d762=\", \"89\u00d738=\"],\n  [\"18\u00d725=\", \"99\u00d776=\"],\n  [\"67\u00d737=\", \"35\u00d755=\"],\n  [\"22\u00d755=\", \"15\u00d794=\"],\n  [\"70\u00d756=\", \"91\u00d799=\"],\n  [\"98\u00d758=\", \"87\u00d798=\"],\n  [\"19\u00d761=\", \"24\u00d723=\"],\n  [\"24\u00d790=\", \"49\u00d766=\"],\n  [\"78\u00d776=\", \"59\u00d760=\"],\n  [\"66\u00d771=\", \"49\u00d752=\"],\n  [\"41\u00d756=\", \"27\u00d724=\"],\n  [\"63\u00d764=\", \"67\u00d734=\"],\n  [\"44\u00d773=\", \"85\u00d735=\"],\n  [\"67\u00d781=\", \"75\u00d724=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each two-digit multiplication problem's operands in the table.\n# Every \"A\u00d7B=\" string occurs exactly once in the document, so a plain\n# Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"15\u00d756=\", \"75\u00d728=\"),\n    @(\"74\u00d737=\", \"53\u00d716=\"),\n    @(\"69\u00d721=\", \"41\u00d771=\"),\n    @(\"19\u00d717=\", \"93\u00d787=\"),\n    @(\"30\u00d778=\", \"55\u00d757=\"),\n    @(\"31\u00d791=\", \"91\u00d768=\"),\n    @(\"40\u00d756=\", \"20\u00d783=\"),\n    @(\"64\u00d753=\", \"57\u00d796=\"),\n    @(\"76\u00d712=\", \"97\u00d716=\"),\n    @(\"50\u00d741=\", \"35\u00d786=\"),\n    @(\"22\u00d732=\", \"89\u00d789=\"),\n    @(\"91\u00d762=\", \"89\u00d738=\"),\n    @(\"18\u00d725=\", \"99\u00d776=\"),\n    @(\"67\u00d737=\", \"35\u00d755=\"),\n    @(\"22\u00d755=\", \"15\u00d794=\"),\n    @(\"70\u00d756=\", \"91\u00d799=\"),\n    @(\"98\u00d758=\", \"87\u00d798=\"),\n    @(\"19\u00d761=\", \"24\u00d723=\"),\n    @(\"24\u00d790=\", \"49\u00d766=\"),\n    @(\"78\u00d776=\", \"59\u00d760=\"),\n    @(\"66\u00d771=\", \"49\u00d752=\"),\n    @(\"41\u00d756=\", \"27\u00d724=\"),\n    @(\"63\u00d764=\", \"67\u00d734=\"),\n    @(\"44\u00d773=\", \"85\u00d735=\"),\n    @(\"67\u00d781=\", \"75\u00d724=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
